# Gantt Chart update: extend/mark "Actual" progress-bar highlighting
# (light-blue fill, same colour as the existing highlighted cells) across
# a number of task rows, plus move the current selection/scroll position
# and bump one row's height.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The highlight colour already used throughout the sheet for "Actual" bars
# (fill FF92CDDC -> BGR 0xDCCD92 for the OLE/VBA color long).
$highlight = 14470546

# Plain fill-only highlights (reuses the sheet's existing "blue fill" style)
$ws.Range("I5").Interior.Color = $highlight
$ws.Range("H7:M7").Interior.Color = $highlight
$ws.Range("H15:M15").Interior.Color = $highlight
$ws.Range("G17").Interior.Color = $highlight
$ws.Range("F19:O19").Interior.Color = $highlight
$ws.Range("I21:K21").Interior.Color = $highlight
$ws.Range("J23:L23").Interior.Color = $highlight
$ws.Range("K25:O25").Interior.Color = $highlight
$ws.Range("M27:N27").Interior.Color = $highlight
$ws.Range("L29:O29").Interior.Color = $highlight
$ws.Range("M31:O31").Interior.Color = $highlight
$ws.Range("M33:O33").Interior.Color = $highlight
$ws.Range("O41").Interior.Color = $highlight

# F5 and O5 get both the highlight fill AND a matching (near-invisible)
# font colour, matching the new combined font+fill style used there.
$ws.Range("F5").Interior.Color = $highlight
$ws.Range("F5").Font.Color = $highlight
$ws.Range("O5").Interior.Color = $highlight
$ws.Range("O5").Font.Color = $highlight

# O27 keeps its existing red font but also gets the highlight fill.
$ws.Range("O27").Interior.Color = $highlight

# Row 23 grew taller.
$ws.Rows.Item(23).RowHeight = 19.5

# Move the view: scroll so row 6 is at the top and select N31.
$ws.Range("N31").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
